# Auto-generated edit script: updates market-price derived columns (H-N)
# across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the
# refreshed Hyperion market data snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 5000
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 5000
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 5000
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -5224
$ws.Range("H14").Value = 5000
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 5000
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 5000
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -5382
$ws.Range("H15").Value = 1104.253
$ws.Range("I15").Value = 1104.253
$ws.Range("K15").Value = 3312.759
$ws.Range("M15").Value = -3143.759
$ws.Range("H18").Value = 1752.9445
$ws.Range("I18").Value = 1838.7059
$ws.Range("K18").Value = 1838.7059
$ws.Range("M18").Value = -1554.7059
$ws.Range("H19").Value = 1723
$ws.Range("I19").Value = 348.5
$ws.Range("K19").Value = 348.5
$ws.Range("M19").Value = -173.5
$ws.Range("H40").Value = 7408.0967
$ws.Range("I40").Value = 21890.2
$ws.Range("J40").Value = 4623.077
$ws.Range("K40").Value = 21890.2
$ws.Range("L40").Value = 4623.077
$ws.Range("M40").Value = -21715.2
$ws.Range("N40").Value = -4973.077
$ws.Range("H51").Value = 6628.32
$ws.Range("J51").Value = 7200.2383
$ws.Range("L51").Value = 7200.2383
$ws.Range("N51").Value = -8168.2383
$ws.Range("H55").Value = 306.9091
$ws.Range("I55").Value = 121.625
$ws.Range("J55").Value = 412.7857
$ws.Range("K55").Value = 121.625
$ws.Range("L55").Value = 412.7857
$ws.Range("M55").Value = 92.375
$ws.Range("N55").Value = -840.7857
$ws.Range("H70").Value = 3001
$ws.Range("I70").Value = 4002
$ws.Range("J70").Value = 2000
$ws.Range("K70").Value = 12006
$ws.Range("L70").Value = 6000
$ws.Range("M70").Value = -11736
$ws.Range("N70").Value = -6540
$ws.Range("H73").Value = 3001
$ws.Range("I73").Value = 4002
$ws.Range("J73").Value = 2000
$ws.Range("K73").Value = 12006
$ws.Range("L73").Value = 6000
$ws.Range("M73").Value = -11070
$ws.Range("N73").Value = -7872
$ws.Range("H76").Value = 5883.278
$ws.Range("I76").Value = 5074.9165
$ws.Range("J76").Value = 7500
$ws.Range("K76").Value = 5074.9165
$ws.Range("L76").Value = 7500
$ws.Range("M76").Value = -4759.9165
$ws.Range("N76").Value = -8130
$ws.Range("H79").Value = 5883.278
$ws.Range("I79").Value = 5074.9165
$ws.Range("J79").Value = 7500
$ws.Range("K79").Value = 5074.9165
$ws.Range("L79").Value = 7500
$ws.Range("M79").Value = -3982.9165
$ws.Range("N79").Value = -9684
$ws.Range("H118").Value = 76923570
$ws.Range("I118").Value = 90909480
$ws.Range("K118").Value = 272728440
$ws.Range("M118").Value = -272726783
$ws.Range("H135").Value = 1020.91895
$ws.Range("I135").Value = 648.6818
$ws.Range("K135").Value = 5838.1362
$ws.Range("M135").Value = -3303.1362
$ws.Range("H138").Value = 3033.125
$ws.Range("J138").Value = 3527.04
$ws.Range("L138").Value = 10581.12
$ws.Range("N138").Value = -20861.12

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4924.47
$ws.Range("I32").Value = 3237.7693
$ws.Range("J32").Value = 21978.889
$ws.Range("K32").Value = 3237.7693
$ws.Range("L32").Value = 21978.889
$ws.Range("M32").Value = -2950.7693
$ws.Range("N32").Value = -22552.889
$ws.Range("H45").Value = 5106
$ws.Range("I45").Value = 2558
$ws.Range("K45").Value = 2558
$ws.Range("M45").Value = -2181
$ws.Range("H97").Value = 12431.632
$ws.Range("I97").Value = 16385.857
$ws.Range("K97").Value = 16385.857
$ws.Range("M97").Value = -15889.857

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 372.8
$ws.Range("I22").Value = 391.25
$ws.Range("J22").Value = 299
$ws.Range("K22").Value = 391.25
$ws.Range("L22").Value = 299
$ws.Range("M22").Value = -218.25
$ws.Range("N22").Value = -645
$ws.Range("H80").Value = 435.84
$ws.Range("J80").Value = 420.23077
$ws.Range("L80").Value = 420.23077
$ws.Range("N80").Value = -2416.23077
$ws.Range("H83").Value = 435.84
$ws.Range("J83").Value = 420.23077
$ws.Range("L83").Value = 2101.15385
$ws.Range("N83").Value = -12085.15385
$ws.Range("H105").Value = 3993.125
$ws.Range("I105").Value = 4324.3335
$ws.Range("K105").Value = 4324.3335
$ws.Range("M105").Value = -2577.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 29410.928
$ws.Range("I31").Value = 23799.4
$ws.Range("J31").Value = 30190.305
$ws.Range("K31").Value = 23799.4
$ws.Range("L31").Value = 30190.305
$ws.Range("M31").Value = -23504.4
$ws.Range("N31").Value = -30780.305
$ws.Range("H34").Value = 29410.928
$ws.Range("I34").Value = 23799.4
$ws.Range("J34").Value = 30190.305
$ws.Range("K34").Value = 23799.4
$ws.Range("L34").Value = 30190.305
$ws.Range("M34").Value = -23597.4
$ws.Range("N34").Value = -30594.305
$ws.Range("H57").Value = 10000
$ws.Range("I57").Value = 10000
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 10000
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -9440
$ws.Range("N57").ClearContents()
$ws.Range("H107").Value = 1565.093
$ws.Range("I107").Value = 1500.1428
$ws.Range("K107").Value = 1500.1428
$ws.Range("M107").Value = 419.8571999999999
$ws.Range("H132").Value = 84963.87
$ws.Range("I132").Value = 68773.87
$ws.Range("J132").Value = 115320.125
$ws.Range("K132").Value = 206321.61
$ws.Range("L132").Value = 345960.375
$ws.Range("M132").Value = -203791.61
$ws.Range("N132").Value = -351020.375
$ws.Range("H134").Value = 7674.353
$ws.Range("I134").Value = 4879.9644
$ws.Range("J134").Value = 20714.834
$ws.Range("K134").Value = 14639.8932
$ws.Range("L134").Value = 62144.50199999999
$ws.Range("M134").Value = -12104.8932
$ws.Range("N134").Value = -67214.50199999999
$ws.Range("H138").Value = 29260
$ws.Range("J138").Value = 29260
$ws.Range("L138").Value = 29260
$ws.Range("N138").Value = -39540

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 737.3333
$ws.Range("J107").Value = 500
$ws.Range("L107").Value = 1500
$ws.Range("N107").Value = -5340

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 292925
$ws.Range("I80").Value = 109812.5
$ws.Range("K80").Value = 109812.5
$ws.Range("M80").Value = -108814.5
$ws.Range("H83").Value = 292925
$ws.Range("I83").Value = 109812.5
$ws.Range("K83").Value = 549062.5
$ws.Range("M83").Value = -544070.5
$ws.Range("H97").Value = 779.1852
$ws.Range("I97").Value = 676.1579
$ws.Range("J97").Value = 1023.875
$ws.Range("K97").Value = 676.1579
$ws.Range("L97").Value = 1023.875
$ws.Range("M97").Value = -180.1579
$ws.Range("N97").Value = -2015.875
$ws.Range("H102").Value = 1542.2174
$ws.Range("I102").Value = 816.71875
$ws.Range("J102").Value = 3200.5
$ws.Range("K102").Value = 816.71875
$ws.Range("L102").Value = 3200.5
$ws.Range("M102").Value = 805.28125
$ws.Range("N102").Value = -6444.5
$ws.Range("H111").Value = 26429
$ws.Range("J111").Value = 26429
$ws.Range("L111").Value = 26429
$ws.Range("N111").Value = -32563
$ws.Range("H122").Value = 289632.78
$ws.Range("I122").Value = 372654.97
$ws.Range("J122").Value = 4985.2856
$ws.Range("K122").Value = 1117964.91
$ws.Range("L122").Value = 14955.8568
$ws.Range("M122").Value = -1115514.91
$ws.Range("N122").Value = -19855.8568
$ws.Range("H132").Value = 8739.868
$ws.Range("J132").Value = 10313.2
$ws.Range("L132").Value = 30939.6
$ws.Range("N132").Value = -35999.60000000001
$ws.Range("H141").Value = 94000
$ws.Range("J141").Value = 94000
$ws.Range("L141").Value = 94000
$ws.Range("N141").Value = -104360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 65489.285
$ws.Range("I22").Value = 178282.6
$ws.Range("J22").Value = 2826.3333
$ws.Range("K22").Value = 178282.6
$ws.Range("L22").Value = 2826.3333
$ws.Range("M22").Value = -177987.6
$ws.Range("N22").Value = -3416.3333
$ws.Range("H27").Value = 65489.285
$ws.Range("I27").Value = 178282.6
$ws.Range("J27").Value = 2826.3333
$ws.Range("K27").Value = 178282.6
$ws.Range("L27").Value = 2826.3333
$ws.Range("M27").Value = -178175.6
$ws.Range("N27").Value = -3040.3333
$ws.Range("H93").Value = 6375.478
$ws.Range("I93").Value = 1935.5625
$ws.Range("K93").Value = 1935.5625
$ws.Range("M93").Value = -687.5625
$ws.Range("H122").Value = 6281.0557
$ws.Range("I122").Value = 4364.3335
$ws.Range("K122").Value = 13093.0005
$ws.Range("M122").Value = -10643.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 839.8182
$ws.Range("I107").Value = 681.05554
$ws.Range("J107").Value = 1030.3334
$ws.Range("K107").Value = 2043.16662
$ws.Range("L107").Value = 3091.0002
$ws.Range("M107").Value = -123.16662
$ws.Range("N107").Value = -6931.0002
$ws.Range("H126").Value = 2828.5
$ws.Range("I126").Value = 2812.9092
$ws.Range("K126").Value = 8438.7276
$ws.Range("M126").Value = -5968.7276
$ws.Range("H132").Value = 203279.84
$ws.Range("I132").Value = 7658.3267
$ws.Range("K132").Value = 22974.9801
$ws.Range("M132").Value = -20444.9801
$ws.Range("H136").Value = 6097.6855
$ws.Range("J136").Value = 3341.1428
$ws.Range("L136").Value = 10023.4284
$ws.Range("N136").Value = -15123.4284
